$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.225
$ws.Range("C2").Value = 55.287
$ws.Range("D2").Value = 1.818
$ws.Range("E2").Value = 1.986
$ws.Range("F2").Value = 76.316

$ws.Range("B3").Value = 0.909
$ws.Range("C3").Value = 1.435
$ws.Range("D3").Value = 0.024
$ws.Range("E3").Value = 0.191
$ws.Range("F3").Value = 2.559

$ws.Range("B4").Value = 4.139
$ws.Range("C4").Value = 9.928000000000001
$ws.Range("D4").Value = 0.67
$ws.Range("E4").Value = 0.167
$ws.Range("F4").Value = 14.904

$ws.Range("B5").Value = 1.818
$ws.Range("C5").Value = 4.019
$ws.Range("D5").Value = 0.12
$ws.Range("E5").Value = 0.263
$ws.Range("F5").Value = 6.220000000000001

$ws.Range("B6").Value = 24.091
$ws.Range("C6").Value = 70.669
$ws.Range("D6").Value = 2.632
$ws.Range("E6").Value = 2.607
$ws.Range("F6").Value = 99.99900000000001
